# Add a new working-hours entry as row 7 (pushing the existing "Total"
# row down to row 8), and update the Total row's SUM formulas to cover
# the extended data range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: relocate the current "Total" row (row 7) down to row 8,
# carrying its values/number-formats/styles along with it. (Using a
# plain Copy instead of Rows.Insert keeps the style table untouched.)
$ws.Range("A7:F7").Copy($ws.Range("A8:F8"))

# Step 2: populate the now-free row 7 with a new time entry, matching
# the look/formatting of the preceding data row (row 6).
$ws.Range("A6:F6").Copy()
$ws.Range("A7:F7").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A7").Value = 45269
$ws.Range("B7").Value = 0.5
$ws.Range("C7").Value = 0.875
$ws.Range("D7").Formula = "=(C7<B7)+C7-B7"
$ws.Range("E7").Value = 10
$ws.Range("F7").Formula = "=(D7*24)*E7"
# Formula entry above can drag in D7's duration number-format; force F7
# back to the plain "Bill" column format it shares with the other rows.
$ws.Range("F7").NumberFormat = "General"

# Step 3: extend the Total row's sums to include the new row 7.
$ws.Range("D8").Formula = "=SUM(D2:D7)"
$ws.Range("F8").Formula = "=SUM(F2:F7)"

$ws.Range("F9").Select() | Out-Null
